$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values in row 1 for the two new columns P and Q, matching the
# existing bordered/bold style used by B1:O1 (style index 1 -> "s=1" cells).
$ws.Range("B1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2-25), flip the I/K/M/O values and populate the two
# new columns P and Q with 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P -> 2
    $ws.Cells.Item($r, 17).Value = 2   # Q -> 2
}
